$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '67.065.79'
$ws.Range('E2').Value = '  +0.72%  '
$ws.Range('D3').Value = '2.531.27'
$ws.Range('E3').Value = '  -2.06%  '
$ws.Range('E4').Value = '  -0.02%  '
$ws.Range('D5').Value = "'589.32"
$ws.Range('E5').Value = '  +0.89%  '
$ws.Range('D6').Value = "'172.97"
$ws.Range('E6').Value = '  +3.86%  '
$ws.Range('E7').Value = '  +0.01%  '
$ws.Range('D8').Value = "'0.527"
$ws.Range('E8').Value = '  +0.11%  '
$ws.Range('D9').Value = '2.531.72'
$ws.Range('E9').Value = '  -2.00%  '
$ws.Range('E10').Value = '  +0.43%  '
$ws.Range('E11').Value = '  +2.14%  '
$ws.Range('D12').Value = "'5.14"
$ws.Range('E12').Value = '  -0.84%  '
$ws.Range('D13').Value = "'0.343"
$ws.Range('E13').Value = '  -3.96%  '
$ws.Range('D14').Value = "'26.53"
$ws.Range('E14').Value = '  -0.85%  '
$ws.Range('D15').Value = '2.992.51'
$ws.Range('E15').Value = '  -2.03%  '
$ws.Range('E16').Value = '  -1.15%  '
$ws.Range('D17').Value = '66.983.75'
$ws.Range('E17').Value = '  +1.00%  '
$ws.Range('D18').Value = '2.519.22'
$ws.Range('E18').Value = '  -2.86%  '
$ws.Range('D19').Value = "'8.12"
$ws.Range('E19').Value = '  +4.65%  '
$ws.Range('D20').Value = "'11.37"
$ws.Range('E20').Value = '  -0.59%  '
$ws.Range('D21').Value = "'354.98"
$ws.Range('E21').Value = '  +0.58%  '
$ws.Range('D22').Value = "'4.18"
$ws.Range('E22').Value = '  -1.35%  '
$ws.Range('D23').Value = "'4.62"
$ws.Range('E23').Value = '  +0.04%  '
$ws.Range('D24').Value = "'1.99"
$ws.Range('E24').Value = '  +5.00%  '
$ws.Range('E25').Value = '  -0.12%  '
$ws.Range('D26').Value = "'69.79"
$ws.Range('E26').Value = '  +1.29%  '
$ws.Range('E27').Value = '  -0.62%  '
$ws.Range('E28').Value = '  -0.07%  '
$ws.Range('D29').Value = '2.655.03'
$ws.Range('E29').Value = '  -2.25%  '
$ws.Range('D30').Value = '0.0₃0976'
$ws.Range('E30').Value = '  -1.21%  '
$ws.Range('D31').Value = "'533.53"
$ws.Range('E31').Value = '  -0.57%  '
$ws.Range('D32').Value = "'8.11"
$ws.Range('E32').Value = '  +1.41%  '
$ws.Range('E33').Value = '  -0.20%  '
$ws.Range('D34').Value = "'1.85"
$ws.Range('E34').Value = '  -0.55%  '
$ws.Range('D35').Value = "'0.131"
$ws.Range('E35').Value = '  -1.02%  '
$ws.Range('D36').Value = "'1.00"
$ws.Range('E36').Value = '  +0.08%  '
$ws.Range('E37').Value = '  -0.50%  '
$ws.Range('D38').Value = "'158.71"
$ws.Range('E38').Value = '  +1.25%  '
$ws.Range('D39').Value = "'18.61"
$ws.Range('E39').Value = '  -0.80%  '
$ws.Range('D40').Value = "'18.45"
$ws.Range('E40').Value = '  +1.13%  '
$ws.Range('D41').Value = "'0.355"
$ws.Range('E41').Value = '  -1.68%  '
$ws.Range('E42').Value = '  +0.17%  '
$ws.Range('D43').Value = "'5.13"
$ws.Range('E43').Value = '  +0.20%  '
$ws.Range('E44').Value = '  -0.06%  '
$ws.Range('D45').Value = "'2.49"
$ws.Range('E45').Value = '  +3.00%  '
$ws.Range('D46').Value = "'149.26"
$ws.Range('E46').Value = '  -0.09%  '
$ws.Range('D47').Value = "'0.556"
$ws.Range('D48').Value = '0.0₆0279'
$ws.Range('E48').Value = '  -2.75%  '
$ws.Range('D49').Value = "'3.68"
$ws.Range('E49').Value = '  -1.25%  '
$ws.Range('D50').Value = "'1.69"
$ws.Range('E50').Value = '  -1.23%  '
$ws.Range('D51').Value = "'0.0757"
$ws.Range('E51').Value = '  -0.25%  '
